$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 40) with the new LeetCode-style entry
$ws.Range("A40").Value = 3469
$ws.Range("B40").Value = "Reschedule Meetings for Maximum Free Time 1"
$ws.Range("C40").Value = "Window Strategy"
$ws.Range("D40").Value = "Collect all gaps[], build a window and slide it to find max free time, we rearrange k meetings merge k+1 gaps, the max size of a window is k+1"

# Update the view scroll/selection state to match the author's saved state
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("D44").Select()
